# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume } for rows whose values changed.
$updates = @{
    2  = @{ D = "27.445.26";     E = "  +2.01%  " }
    3  = @{ D = "1.835.80";      E = "  +1.12%  " }
    4  = @{               E = "  +1.02%  " }
    5  = @{ D = "314.35";        E = "  +1.73%  " }
    6  = @{               E = "  +0.91%  " }
    7  = @{               E = "  +2.02%  " }
    8  = @{ D = "0.3691";        E = "  +0.76%  " }
    9  = @{ D = "0.07457";       E = "  +1.40%  " }
    10 = @{ D = "0.8854";        E = "  +1.84%  " }
    11 = @{               E = "  +0.90%  " }
    12 = @{ D = "1.889.16";      E = "  +0.76%  " }
    13 = @{ D = "0.07334";       E = "  +3.41%  " }
    14 = @{ D = "5.453";         E = "  +1.43%  " }
    15 = @{ D = "93.17";         E = "  +1.87%  " }
    16 = @{ D = "6.582" }
    17 = @{ D = "1.011";         E = "  +0.74%  " }
    18 = @{ D = "0.000008822";   E = "  +1.27%  " }
    19 = @{               E = "  +0.93%  " }
    20 = @{ D = "27.672.56";     E = "  +2.67%  " }
    21 = @{ D = "14.81";         E = "  +1.10%  " }
    22 = @{ D = "5.318";         E = "  +0.32%  " }
    23 = @{               E = "  +0.58%  " }
    24 = @{ D = "2.103.62";      E = "  +1.41%  " }
    25 = @{ D = "1.912";         E = "  +0.94%  " }
    26 = @{ D = "152.13";        E = "  +0.86%  " }
    27 = @{ D = "18.65";         E = "  +1.82%  " }
    28 = @{               E = "  +0.58%  " }
    29 = @{ D = "5.245";         E = "  -0.41%  " }
    30 = @{ D = "117.47";        E = "  +2.08%  " }
    31 = @{ D = "0.08998";       E = "  +1.13%  " }
    32 = @{ D = "0.7568";        E = "  -0.07%  " }
    33 = @{               E = "  +1.80%  " }
    34 = @{ D = "4.553";         E = "  +1.51%  " }
    35 = @{ D = "2.947";         E = "  +1.26%  " }
    36 = @{               E = "  +1.06%  " }
    37 = @{ D = "1.104" }
    38 = @{ D = "0.05345";       E = "  +1.18%  " }
    39 = @{ D = "0.01954";       E = "  +0.35%  " }
    40 = @{ D = "2.980";         E = "  -0.19%  " }
    41 = @{ D = "7.330";         E = "  +1.26%  " }
    42 = @{ D = "2.409";         E = "  +5.30%  " }
    43 = @{ D = "0.5332";        E = "  +0.64%  " }
    44 = @{ D = "0.1660";        E = "  +0.33%  " }
    45 = @{ D = "8.519";         E = "  +1.02%  " }
    46 = @{ D = "0.4910";        E = "  +0.88%  " }
    47 = @{ D = "10.60";         E = "  +1.29%  " }
    48 = @{               E = "  +1.05%  " }
    49 = @{ D = "104.98";        E = "  +1.59%  " }
    50 = @{ D = "1.677";         E = "  +1.00%  " }
    51 = @{ D = "0.06303";       E = "  +0.16%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        # Force a text number format first so Excel doesn't reinterpret the
        # price string (e.g. "314.35") as a numeric value; then restore the
        # cell's style to Normal so no stray formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["D"]
        $cell.Style = "Normal"
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Range("E$row").Value = $rowData["E"]
    }
}
